$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Values shared across the new rows (columns B through J), identical to row 106
$rowValues = @(116.4121952, 0.00170247, 0.008850780000000001, 0.06933635, 12792.90181321, 465.80531254, 0.24, 1.7904431, 485.38834923)

# Starting date serial for new row 107 is 45663, incrementing by 1 through row 114 (45670)
$startSerial = 45663
$startRow = 107
$endRow = 114

for ($r = $startRow; $r -le $endRow; $r++) {
    $serial = $startSerial + ($r - $startRow)

    # Column A: copy formatting (date number format/style) from the row above, then set the value
    $ws.Cells.Item($r - 1, 1).Copy()
    $ws.Cells.Item($r, 1).PasteSpecial(-4122)
    $ws.Cells.Item($r, 1).Value = $serial

    # Columns B through J
    for ($c = 2; $c -le 10; $c++) {
        $ws.Cells.Item($r, $c).Value = $rowValues[$c - 2]
    }
}

$excel.CutCopyMode = 0
